$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(602).Copy()
$ws.Rows.Item(603).PasteSpecial(-4122)
$ws.Range("A603").Value = "2021-11-29"
$ws.Range("B603").Value = 74
$ws.Range("C603").Value = 73
